# Update column F ("dSF") values for the specified rows to reflect
# repulled data / recalculated mean values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 5
    13 = -5
    15 = -3
    21 = 0
    26 = -1
    27 = -1
    32 = 2
    35 = 4
    36 = -2
    37 = 0
    39 = 0
    41 = -1
    42 = -8
    46 = -1
    52 = 4
    56 = 2
    59 = -2
    63 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
